$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix translation typo in the "Other Unit Of Measure" header (cell H1)
$ws.Range("H1").Value = "Other Unit Of Measure وحدة قياس أخرى"

# 2. Update the sheet view: scroll/selection moves from F2 (topLeftCell B1) to G1
$ws.Range("G1").Select()

# 3. Add a new blank data validation rule for the whole column C (C1:C1048576).
#    This also causes the original rule that covered "C2:F5 H2:I2 G1 K1" to lose
#    its C2:C5 portion, leaving it covering "D2:F5 H2:I2 G1 K1".
$colCValidation = $ws.Range("C1:C1048576").Validation
$colCValidation.Delete()
$colCValidation.Add(0, 1, 1)
$colCValidation.IgnoreBlank = $true
$colCValidation.InCellDropdown = $true
$colCValidation.ShowInput = $true
$colCValidation.ShowError = $true
